# The upstream MWW (Mann-Whitney/Wilcoxon) test results were recomputed and the
# per-nutrient blocks of pairwise cluster-comparison rows (14-103, six rows each:
# comparisons 1-2, 1-3, 1-4, 2-3, 2-4, 3-4) were re-assigned to a new nutrient
# ordering (Na, K, Mg, Fe, Zn, Cu, Mn, As, Se, Ni, Co, Sr, Cd, Pb, Ag). Column A
# (nutrient label), D (alpha_MW p-value) and E (significant yes/no) are updated
# per block; columns B and C (the cluster-pair ids) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: nutrient name, then the six [alpha_MW, significant] pairs for the
# comparisons 1-2, 1-3, 1-4, 2-3, 2-4, 3-4 (in that row order).
$blocks = @(
    @{ Name = "Na"; Stats = @(@(0.00910269640919796, "yes"), @(0.0145442792501616, "yes"), @(0.0000833531793284115, "yes"), @(0.851814851814852, "no"), @(0.000329082682023858, "yes"), @(0.035964035964036, "yes")) }
    @{ Name = "K"; Stats = @(@(0.000793839803127729, "yes"), @(0.462346477052359, "no"), @(0.710300071445582, "no"), @(0.00266400266400266, "yes"), @(0.00156314273961333, "yes"), @(0.388411588411588, "no")) }
    @{ Name = "Mg"; Stats = @(@(0.237543330422587, "no"), @(0.732546864899806, "no"), @(0.0000238151940938319, "yes"), @(0.754578754578755, "no"), @(0.000164541341011929, "yes"), @(0.000799200799200799, "yes")) }
    @{ Name = "Fe"; Stats = @(@(0.395041147363129, "no"), @(0.807530704589528, "no"), @(0.000797809002143367, "yes"), @(0.572760572760573, "no"), @(0.00551213492389963, "yes"), @(0.0015984015984016, "yes")) }
    @{ Name = "Zn"; Stats = @(@(0.000105845307083697, "yes"), @(0.0782159017453135, "no"), @(0.000226244343891403, "yes"), @(0.228438228438228, "no"), @(0.276593994241053, "no"), @(0.327672327672328, "no")) }
    @{ Name = "Cu"; Stats = @(@(0.544759334233018, "no"), @(0.180187459599224, "no"), @(0.00115503691355085, "yes"), @(0.572760572760573, "no"), @(0.0273961332784862, "yes"), @(0.223776223776224, "no")) }
    @{ Name = "Mn"; Stats = @(@(0.151755709031251, "no"), @(0.660471881060116, "no"), @(0.00228625863300786, "yes"), @(0.413586413586414, "no"), @(0.480707527766351, "no"), @(0.0027972027972028, "yes")) }
    @{ Name = "As"; Stats = @(@(0.0000264613267709243, "yes"), @(0.0476729153199741, "yes"), @(0.00565610859728507, "yes"), @(0.000666000666000666, "yes"), @(0.0000822706705059646, "yes"), @(0.0003996003996004, "yes")) }
    @{ Name = "Se"; Stats = @(@(0.206001428911646, "no"), @(0.807530704589528, "no"), @(0.710300071445582, "no"), @(0.0592740592740593, "no"), @(0.0927190456602221, "no"), @(0.863936063936064, "no")) }
    @{ Name = "Ni"; Stats = @(@(0.457193123368074, "no"), @(0.000323206205559147, "yes"), @(0.766432483924744, "no"), @(0.00133200133200133, "yes"), @(0.541423282599753, "no"), @(0.035964035964036, "yes")) }
    @{ Name = "Co"; Stats = @(@(0.341689072477986, "no"), @(0.00108166488475667, "yes"), @(0.361573522385796, "no"), @(0.000666000666000666, "yes"), @(1, "no"), @(0.0290951317901257, "yes")) }
    @{ Name = "Sr"; Stats = @(@(0.0203487602868408, "yes"), @(0.660471881060116, "no"), @(0.130804953560372, "no"), @(0.0126540126540127, "yes"), @(0.235870012340601, "no"), @(0.0663336663336663, "no")) }
    @{ Name = "Cd"; Stats = @(@(0.363511807596544, "no"), @(0.0273109243697479, "yes"), @(0.111943319838057, "no"), @(0.228438228438228, "no"), @(0.0152200740436035, "yes"), @(0.0027972027972028, "yes")) }
    @{ Name = "Pb"; Stats = @(@(0.00383653803899588, "yes"), @(0.404330963154493, "no"), @(0.196343631890015, "no"), @(0.0592740592740593, "no"), @(0.00810146447141977, "yes"), @(0.215514712957383, "no")) }
    @{ Name = "Ag"; Stats = @(@(0.535186298659275, "no"), @(0.0498760230060749, "yes"), @(0.361392678921389, "no"), @(0.01998001998002, "yes"), @(0.177665829328152, "no"), @(0.0495504495504496, "yes")) }
)

$row = 14
foreach ($block in $blocks) {
    foreach ($stat in $block.Stats) {
        $ws.Range("A$row").Value = $block.Name
        $ws.Range("D$row").Value = $stat[0]
        $ws.Range("E$row").Value = $stat[1]
        $row = $row + 1
    }
}

Write-Host "Updated rows 14 through $($row - 1)"
